$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "$ 27.386 CLP 29-10-20"
$ws.Range("A48").Value = "$ 27.391 CLP 30-10-20"
$ws.Range("A49").Value = "$ 27.391 CLP 30-10-20"
